$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-12 (row 25)
$ws.Range("B25").Value = 6463
$ws.Range("D25").Value = 6009155
$ws.Range("E25").Value = 929.7779668884419
$ws.Range("F25").Value = 9.709726701748433
$ws.Range("H25").Value = 25.85023684533492
